$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text (matches original inlineStr type)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.093.00"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.86"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.92"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.99"
$ws.Range("E8").Value = "  +10.02%  "
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.125.68"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.74"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.46"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.682"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.037.64"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.42"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.68"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.77"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.59"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  +23.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.73"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("E29").Value = "  +2.77%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.02"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.01"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.04"
$ws.Range("E34").Value = "  +14.10%  "
$ws.Range("E35").Value = "  +23.49%  "
$ws.Range("E36").Value = "  +12.77%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.09"
$ws.Range("E38").Value = "  +14.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.83"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.352.39"
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.12"
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("E43").Value = "  +6.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.86"
$ws.Range("E44").Value = "  +62.10%  "
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0543"
$ws.Range("E47").Value = "  +5.98%  "
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.036.06"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.45"
$ws.Range("E51").Value = "  +16.35%  "
